$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("L6").Value = 5.5
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 7
$ws.Range("Q6").Value = 2.3
$ws.Range("R6").Value = 1.6
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("AC6").Value = 7
$ws.Range("AE6").Value = 19
$ws.Range("AR6").Value = 67
$ws.Range("AW6").Value = 6.5
$ws.Range("AX6").Value = 29
$ws.Range("BA6").Value = 151

# Row 15
$ws.Range("M15").Value = 1.03
$ws.Range("N15").Value = 17
$ws.Range("O15").Value = 1.17
$ws.Range("P15").Value = 5
$ws.Range("Q15").Value = 1.57
$ws.Range("R15").Value = 2.35

# Row 16
$ws.Range("G16").Value = 2.35
$ws.Range("H16").Value = 3.6
$ws.Range("I16").Value = 2.7
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 2.3
$ws.Range("L16").Value = 3.2
$ws.Range("N16").Value = 15
$ws.Range("O16").Value = 1.18
$ws.Range("P16").Value = 4.5
$ws.Range("Q16").Value = 1.62
$ws.Range("R16").Value = 2.25
$ws.Range("S16").Value = 1.3
$ws.Range("T16").Value = 3.4
$ws.Range("U16").Value = 1.53
$ws.Range("V16").Value = 2.38
$ws.Range("W16").Value = 11
$ws.Range("Z16").Value = 23
$ws.Range("AA16").Value = 17
$ws.Range("AB16").Value = 23
$ws.Range("AC16").Value = 15
$ws.Range("AD16").Value = 7
$ws.Range("AG16").Value = 126
$ws.Range("AH16").Value = 12
$ws.Range("AI16").Value = 15
$ws.Range("AJ16").Value = 10
$ws.Range("AL16").Value = 19
$ws.Range("AO16").Value = 13
$ws.Range("AP16").Value = 19
$ws.Range("AQ16").Value = 41
$ws.Range("AR16").Value = 51
$ws.Range("AT16").Value = 3.4
$ws.Range("AW16").Value = 5
$ws.Range("AY16").Value = 21
$ws.Range("BC16").Value = 351

# Row 17
$ws.Range("N17").Value = 12

# Row 21
$ws.Range("M21").Value = 1.03
$ws.Range("N21").Value = 17
$ws.Range("Q21").Value = 1.53
$ws.Range("R21").Value = 2.4

# Row 22
$ws.Range("G22").Value = 1.33
$ws.Range("H22").Value = 5.5
$ws.Range("I22").Value = 7
$ws.Range("N22").Value = 23
$ws.Range("Q22").Value = 1.36
$ws.Range("R22").Value = 3.1
$ws.Range("X22").Value = 9
$ws.Range("AC22").Value = 23
$ws.Range("AH22").Value = 26
$ws.Range("AO22").Value = 6.5
$ws.Range("AW22").Value = 9

# Row 23
$ws.Range("O23").Value = 1.4
$ws.Range("P23").Value = 2.75
$ws.Range("Q23").Value = 2.25
$ws.Range("R23").Value = 1.62

# Row 31
$ws.Range("G31").Value = 2.3
$ws.Range("H31").Value = 3.7
$ws.Range("I31").Value = 2.8
$ws.Range("U31").Value = 1.44
$ws.Range("V31").Value = 2.63
$ws.Range("AC31").Value = 19
$ws.Range("AM31").Value = 23

# Row 32
$ws.Range("M32").Value = 1.08
$ws.Range("N32").Value = 8
$ws.Range("Q32").Value = 2.25
$ws.Range("R32").Value = 1.62

# Row 33
$ws.Range("O33").Value = 1.22
$ws.Range("P33").Value = 4
$ws.Range("Q33").Value = 1.7
$ws.Range("R33").Value = 2.1

# Row 35
$ws.Range("J35").Value = 2.87

# Row 36
$ws.Range("G36").Value = 2.05
$ws.Range("H36").Value = 3.2
$ws.Range("I36").Value = 3.45
$ws.Range("J36").Value = 2.6
$ws.Range("K36").Value = 2.07
$ws.Range("L36").Value = 3.95
$ws.Range("N36").Value = 9.8
$ws.Range("O36").Value = 1.33
$ws.Range("P36").Value = 2.8
$ws.Range("Q36").Value = 1.98
$ws.Range("R36").Value = 1.65
$ws.Range("S36").Value = 1.39
$ws.Range("T36").Value = 2.55
$ws.Range("V36").Value = 1.83
$ws.Range("W36").Value = 7
$ws.Range("X36").Value = 9.75
$ws.Range("Z36").Value = 19
$ws.Range("AA36").Value = 17
$ws.Range("AD36").Value = 6.2
$ws.Range("AE36").Value = 14.5
$ws.Range("AH36").Value = 9.5
$ws.Range("AI36").Value = 18
$ws.Range("AJ36").Value = 12
$ws.Range("AK36").Value = 50
$ws.Range("AM36").Value = 40
$ws.Range("AN36").Value = 3.9
$ws.Range("AO36").Value = 10.25
$ws.Range("AP36").Value = 18.5
$ws.Range("AQ36").Value = 40
$ws.Range("AR36").Value = 70
$ws.Range("AT36").Value = 2.5
$ws.Range("AW36").Value = 5.3
$ws.Range("AX36").Value = 19.5
$ws.Range("AZ36").Value = 100
